$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @(
    @{Row=11; I='aa'; J='Agree/Accept'},
    @{Row=16; I='aa'; J='Agree/Accept'},
    @{Row=24; I='aa'; J='Agree/Accept'},
    @{Row=40; I='aa'; J='Agree/Accept'},
    @{Row=43; I='aa'; J='Agree/Accept'},
    @{Row=45; I='aa'; J='Agree/Accept'},
    @{Row=57; I='sd'; J='Statement-non-opinion'},
    @{Row=69; I='sv'; J='Statement-opinion'},
    @{Row=102; I='aa'; J='Agree/Accept'},
    @{Row=106; I='aa'; J='Agree/Accept'},
    @{Row=111; I='%'; J='Uninterpretable'},
    @{Row=133; I='aa'; J='Agree/Accept'},
    @{Row=135; I='aa'; J='Agree/Accept'},
    @{Row=150; I='aa'; J='Agree/Accept'},
    @{Row=152; I='aa'; J='Agree/Accept'},
    @{Row=153; I='sd'; J='Statement-non-opinion'},
    @{Row=157; I='aa'; J='Agree/Accept'},
    @{Row=161; I='aa'; J='Agree/Accept'},
    @{Row=164; I='aa'; J='Agree/Accept'},
    @{Row=177; I='sd'; J='Statement-non-opinion'},
    @{Row=178; I='aa'; J='Agree/Accept'},
    @{Row=183; I='b'; J='Acknowledge (Backchannel)'},
    @{Row=185; I='aa'; J='Agree/Accept'},
    @{Row=190; I='sd'; J='Statement-non-opinion'},
    @{Row=196; I='sv'; J='Statement-opinion'},
    @{Row=197; I='sv'; J='Statement-opinion'},
    @{Row=205; I='sv'; J='Statement-opinion'},
    @{Row=213; I='sd'; J='Statement-non-opinion'},
    @{Row=217; I='aa'; J='Agree/Accept'},
    @{Row=239; I='sd'; J='Statement-non-opinion'},
    @{Row=240; I='sd'; J='Statement-non-opinion'},
    @{Row=249; I='b'; J='Acknowledge (Backchannel)'},
    @{Row=261; I='%'; J='Uninterpretable'},
    @{Row=265; I='sv'; J='Statement-opinion'},
    @{Row=268; I='sv'; J='Statement-opinion'},
    @{Row=271; I='ba'; J='Appreciation'},
    @{Row=283; I='sv'; J='Statement-opinion'},
    @{Row=285; I='sv'; J='Statement-opinion'},
    @{Row=295; I='sv'; J='Statement-opinion'},
    @{Row=313; I='ba'; J='Appreciation'},
    @{Row=321; I='sd'; J='Statement-non-opinion'},
    @{Row=322; I='ba'; J='Appreciation'},
    @{Row=326; I='aa'; J='Agree/Accept'},
    @{Row=330; I='aa'; J='Agree/Accept'},
    @{Row=339; I='%'; J='Uninterpretable'},
    @{Row=342; I='aa'; J='Agree/Accept'},
    @{Row=346; I='aa'; J='Agree/Accept'},
    @{Row=353; I='aa'; J='Agree/Accept'},
    @{Row=360; I='aa'; J='Agree/Accept'},
    @{Row=371; I='%'; J='Uninterpretable'},
    @{Row=377; I='sd'; J='Statement-non-opinion'},
    @{Row=379; I='sd'; J='Statement-non-opinion'},
    @{Row=398; I='sv'; J='Statement-opinion'},
    @{Row=399; I='aa'; J='Agree/Accept'},
    @{Row=404; I='aa'; J='Agree/Accept'},
    @{Row=411; I='sd'; J='Statement-non-opinion'},
    @{Row=413; I='aa'; J='Agree/Accept'},
    @{Row=419; I='aa'; J='Agree/Accept'},
    @{Row=425; I='sd'; J='Statement-non-opinion'},
    @{Row=427; I='sv'; J='Statement-opinion'},
    @{Row=429; I='sv'; J='Statement-opinion'},
    @{Row=434; I='sv'; J='Statement-opinion'}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.I
    $ws.Cells.Item($u.Row, 10).Value = $u.J
}

Write-Output "Updated $($updates.Count) rows"
